$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 27.1666666666667
$ws.Columns.Item(2).ColumnWidth = 21.1666666666667
$ws.Columns.Item(3).ColumnWidth = 20.1666666666667
$ws.Columns.Item(4).ColumnWidth = 20.1666666666667
$ws.Columns.Item(5).ColumnWidth = 19.1666666666667
$ws.Columns.Item(6).ColumnWidth = 25.1666666666667
$ws.Columns.Item(7).ColumnWidth = 32.1666666666667
$ws.Columns.Item(8).ColumnWidth = 26.1666666666667
$ws.Columns.Item(9).ColumnWidth = 41.1666666666667

# --- Clear old content that is being displaced (old GBPUSD/USDJPY/XAUUSD blocks at rows 8-28) ---
$ws.Range("A8:I28").Clear()

# --- Materialize empty placeholder rows across the whole used range (outline level 0), matching source export ---
$ws.Range("A8:A118").EntireRow.OutlineLevel = 0

# ===== Block: EURUSD =====
$ws.Cells.Item(1,1).Value2 = "EURUSD"
$ws.Cells.Item(2,2).Value2 = "count_of_occurrences"
$ws.Cells.Item(2,3).Value2 = "PnL_per_lot"
$ws.Cells.Item(2,4).Value2 = "total_profit"
$ws.Cells.Item(2,5).Value2 = "total_volume"
$ws.Cells.Item(2,6).Value2 = "typical_spread_in_points"
$ws.Cells.Item(2,7).Value2 = "weighted_avg_execution_spread_`$"
$ws.Cells.Item(2,8).Value2 = "percentage_of_occurrences"
$ws.Cells.Item(2,9).Value2 = "one_point_increase_of_weighted_spread_lr"
$ws.Cells.Item(3,1).Value2 = "Volatility_Trend"
$ws.Cells.Item(4,1).Value2 = "High Volatility + No Trend"
$ws.Cells.Item(4,2).Value2 = 197
$ws.Cells.Item(4,3).Value2 = -8.288711069913706
$ws.Cells.Item(4,4).Value2 = -34979435.88
$ws.Cells.Item(4,5).Value2 = 4047717.57
$ws.Cells.Item(4,6).Value2 = 9.145454545452539
$ws.Cells.Item(4,7).Value2 = 8.850136861443799
$ws.Cells.Item(4,8).Value2 = 35.11586452762923
$ws.Cells.Item(4,9).Value2 = -1.479759821313245
$ws.Cells.Item(5,1).Value2 = "High Volatility + Trend"
$ws.Cells.Item(5,2).Value2 = 108
$ws.Cells.Item(5,3).Value2 = -21.49854362216667
$ws.Cells.Item(5,4).Value2 = -50969123.26
$ws.Cells.Item(5,5).Value2 = 2243046.35
$ws.Cells.Item(5,6).Value2 = 9.225806451613174
$ws.Cells.Item(5,7).Value2 = 8.826627322132296
$ws.Cells.Item(5,8).Value2 = 19.25133689839572
$ws.Cells.Item(5,9).Value2 = -1.659191439800495
$ws.Cells.Item(6,1).Value2 = "Low Volatility + No Trend"
$ws.Cells.Item(6,2).Value2 = 156
$ws.Cells.Item(6,3).Value2 = 2.362199022230769
$ws.Cells.Item(6,4).Value2 = 6939455.84
$ws.Cells.Item(6,5).Value2 = 2476737.09
$ws.Cells.Item(6,6).Value2 = 11.01724137930831
$ws.Cells.Item(6,7).Value2 = 9.375202851878555
$ws.Cells.Item(6,8).Value2 = 27.80748663101604
$ws.Cells.Item(6,9).Value2 = 0.09125545721933777
$ws.Cells.Item(7,1).Value2 = "Low Volatility + Trend"
$ws.Cells.Item(7,2).Value2 = 100
$ws.Cells.Item(7,3).Value2 = -5.37864698262
$ws.Cells.Item(7,4).Value2 = -7614650.48
$ws.Cells.Item(7,5).Value2 = 1375866.59
$ws.Cells.Item(7,6).Value2 = 9.899999999999974
$ws.Cells.Item(7,7).Value2 = 8.457566604326196
$ws.Cells.Item(7,8).Value2 = 17.825311942959
$ws.Cells.Item(7,9).Value2 = -1.876264873612548

# ===== Block: GBPUSD =====
$ws.Range("A1").Copy()
$ws.Cells.Item(38,1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(38,1).Value2 = "GBPUSD"
$ws.Cells.Item(39,2).Value2 = "count_of_occurrences"
$ws.Cells.Item(39,3).Value2 = "PnL_per_lot"
$ws.Cells.Item(39,4).Value2 = "total_profit"
$ws.Cells.Item(39,5).Value2 = "total_volume"
$ws.Cells.Item(39,6).Value2 = "typical_spread_in_points"
$ws.Cells.Item(39,7).Value2 = "weighted_avg_execution_spread_`$"
$ws.Cells.Item(39,8).Value2 = "percentage_of_occurrences"
$ws.Cells.Item(39,9).Value2 = "one_point_increase_of_weighted_spread_lr"
$ws.Cells.Item(40,1).Value2 = "Volatility_Trend"
$ws.Cells.Item(41,1).Value2 = "High Volatility + No Trend"
$ws.Cells.Item(41,2).Value2 = 214
$ws.Cells.Item(41,3).Value2 = -4.070678911214953
$ws.Cells.Item(41,4).Value2 = -12741302.841987
$ws.Cells.Item(41,5).Value2 = 2599682.23
$ws.Cells.Item(41,6).Value2 = 10.43617021276468
$ws.Cells.Item(41,7).Value2 = 12.04392561918357
$ws.Cells.Item(41,8).Value2 = 38.14616755793227
$ws.Cells.Item(41,9).Value2 = -0.02839638508320652
$ws.Cells.Item(42,1).Value2 = "High Volatility + Trend"
$ws.Cells.Item(42,2).Value2 = 82
$ws.Cells.Item(42,3).Value2 = -44.0579731097561
$ws.Cells.Item(42,4).Value2 = -48833532.44791
$ws.Cells.Item(42,5).Value2 = 947369.55
$ws.Cells.Item(42,6).Value2 = 12.45454545454414
$ws.Cells.Item(42,7).Value2 = 12.58185417326494
$ws.Cells.Item(42,8).Value2 = 14.61675579322638
$ws.Cells.Item(42,9).Value2 = -3.488446647930592
$ws.Cells.Item(43,1).Value2 = "Low Volatility + No Trend"
$ws.Cells.Item(43,2).Value2 = 234
$ws.Cells.Item(43,3).Value2 = -1.380881414529914
$ws.Cells.Item(43,4).Value2 = -3001446.303141001
$ws.Cells.Item(43,5).Value2 = 2448729.54
$ws.Cells.Item(43,6).Value2 = 13.47222222222156
$ws.Cells.Item(43,7).Value2 = 13.47261067135606
$ws.Cells.Item(43,8).Value2 = 41.71122994652406
$ws.Cells.Item(43,9).Value2 = -0.006549058777956531
$ws.Cells.Item(44,1).Value2 = "Low Volatility + Trend"
$ws.Cells.Item(44,2).Value2 = 31
$ws.Cells.Item(44,3).Value2 = -42.58803487096774
$ws.Cells.Item(44,4).Value2 = -17835065.53399
$ws.Cells.Item(44,5).Value2 = 323749.73
$ws.Cells.Item(44,6).Value2 = 11
$ws.Cells.Item(44,7).Value2 = 15.6018364113038
$ws.Cells.Item(44,8).Value2 = 5.525846702317291
$ws.Cells.Item(44,9).Value2 = -1.016898904927557

# ===== Block: USDJPY =====
$ws.Range("A1").Copy()
$ws.Cells.Item(75,1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(75,1).Value2 = "USDJPY"
$ws.Cells.Item(76,2).Value2 = "count_of_occurrences"
$ws.Cells.Item(76,3).Value2 = "PnL_per_lot"
$ws.Cells.Item(76,4).Value2 = "total_profit"
$ws.Cells.Item(76,5).Value2 = "total_volume"
$ws.Cells.Item(76,6).Value2 = "typical_spread_in_points"
$ws.Cells.Item(76,7).Value2 = "weighted_avg_execution_spread_`$"
$ws.Cells.Item(76,8).Value2 = "percentage_of_occurrences"
$ws.Cells.Item(76,9).Value2 = "one_point_increase_of_weighted_spread_lr"
$ws.Cells.Item(77,1).Value2 = "Volatility_Trend"
$ws.Cells.Item(78,1).Value2 = "High Volatility + No Trend"
$ws.Cells.Item(78,2).Value2 = 181
$ws.Cells.Item(78,3).Value2 = -13.67256108839779
$ws.Cells.Item(78,4).Value2 = -13317215
$ws.Cells.Item(78,5).Value2 = 901024.2
$ws.Cells.Item(78,6).Value2 = 13.17333333333175
$ws.Cells.Item(78,7).Value2 = 10.61131960572962
$ws.Cells.Item(78,8).Value2 = 32.26381461675579
$ws.Cells.Item(78,9).Value2 = -4.045946811637981
$ws.Cells.Item(79,1).Value2 = "High Volatility + Trend"
$ws.Cells.Item(79,2).Value2 = 103
$ws.Cells.Item(79,3).Value2 = -35.49163069902912
$ws.Cells.Item(79,4).Value2 = -22044239
$ws.Cells.Item(79,5).Value2 = 635193.35
$ws.Cells.Item(79,6).Value2 = 12.72222222222133
$ws.Cells.Item(79,7).Value2 = 9.56386993763122
$ws.Cells.Item(79,8).Value2 = 18.36007130124777
$ws.Cells.Item(79,9).Value2 = -2.353777940763052
$ws.Cells.Item(80,1).Value2 = "Low Volatility + No Trend"
$ws.Cells.Item(80,2).Value2 = 188
$ws.Cells.Item(80,3).Value2 = -2.907438643617021
$ws.Cells.Item(80,4).Value2 = -1269776
$ws.Cells.Item(80,5).Value2 = 418436.7
$ws.Cells.Item(80,6).Value2 = 15.89285714285684
$ws.Cells.Item(80,7).Value2 = 10.40518699754833
$ws.Cells.Item(80,8).Value2 = 33.51158645276293
$ws.Cells.Item(80,9).Value2 = -1.075074353199468
$ws.Cells.Item(81,1).Value2 = "Low Volatility + Trend"
$ws.Cells.Item(81,2).Value2 = 89
$ws.Cells.Item(81,3).Value2 = -20.26355995505618
$ws.Cells.Item(81,4).Value2 = -4740458
$ws.Cells.Item(81,5).Value2 = 256229.87
$ws.Cells.Item(81,6).Value2 = 11.93333333332963
$ws.Cells.Item(81,7).Value2 = 9.942112108275058
$ws.Cells.Item(81,8).Value2 = 15.86452762923351
$ws.Cells.Item(81,9).Value2 = -6.200267236059804

# ===== Block: XAUUSD =====
$ws.Range("A1").Copy()
$ws.Cells.Item(112,1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(112,1).Value2 = "XAUUSD"
$ws.Cells.Item(113,2).Value2 = "count_of_occurrences"
$ws.Cells.Item(113,3).Value2 = "PnL_per_lot"
$ws.Cells.Item(113,4).Value2 = "total_profit"
$ws.Cells.Item(113,5).Value2 = "total_volume"
$ws.Cells.Item(113,6).Value2 = "typical_spread_in_points"
$ws.Cells.Item(113,7).Value2 = "weighted_avg_execution_spread_`$"
$ws.Cells.Item(113,8).Value2 = "percentage_of_occurrences"
$ws.Cells.Item(113,9).Value2 = "one_point_increase_of_weighted_spread_lr"
$ws.Cells.Item(114,1).Value2 = "Volatility_Trend"
$ws.Cells.Item(115,1).Value2 = "High Volatility + No Trend"
$ws.Cells.Item(115,2).Value2 = 223
$ws.Cells.Item(115,3).Value2 = -30.3349935470852
$ws.Cells.Item(115,4).Value2 = -167811765.68
$ws.Cells.Item(115,5).Value2 = 4905503.47
$ws.Cells.Item(115,6).Value2 = 20.84745762712242
$ws.Cells.Item(115,7).Value2 = 27.78015214876156
$ws.Cells.Item(115,8).Value2 = 40.0359066427289
$ws.Cells.Item(115,9).Value2 = 0.174580775022823
$ws.Cells.Item(116,1).Value2 = "High Volatility + Trend"
$ws.Cells.Item(116,2).Value2 = 141
$ws.Cells.Item(116,3).Value2 = -56.62320707801418
$ws.Cells.Item(116,4).Value2 = -201042079.38
$ws.Cells.Item(116,5).Value2 = 3039428.76
$ws.Cells.Item(116,6).Value2 = 20.94736842105631
$ws.Cells.Item(116,7).Value2 = 31.46887424561551
$ws.Cells.Item(116,8).Value2 = 25.31418312387791
$ws.Cells.Item(116,9).Value2 = -0.8735817132808531
$ws.Cells.Item(117,1).Value2 = "Low Volatility + No Trend"
$ws.Cells.Item(117,2).Value2 = 132
$ws.Cells.Item(117,3).Value2 = 14.67839834848485
$ws.Cells.Item(117,4).Value2 = 44661475.96
$ws.Cells.Item(117,5).Value2 = 2853821.63
$ws.Cells.Item(117,6).Value2 = 21.87500000000284
$ws.Cells.Item(117,7).Value2 = 26.55715673764307
$ws.Cells.Item(117,8).Value2 = 23.6983842010772
$ws.Cells.Item(117,9).Value2 = 0.0995374841512332
$ws.Cells.Item(118,1).Value2 = "Low Volatility + Trend"
$ws.Cells.Item(118,2).Value2 = 61
$ws.Cells.Item(118,3).Value2 = -0.7356025245901638
$ws.Cells.Item(118,4).Value2 = 2137769.51
$ws.Cells.Item(118,5).Value2 = 1340254.44
$ws.Cells.Item(118,6).Value2 = 22.15384615384874
$ws.Cells.Item(118,7).Value2 = 26.91615958414931
$ws.Cells.Item(118,8).Value2 = 10.95152603231598
$ws.Cells.Item(118,9).Value2 = -1.499930589854987

$ws.Range("A1").Select()